$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = -517866088

$ws.Range("C18").Formula = "=SUM(C12:C17)"

$ws.Range("C19").Value = -392700000

$ws.Range("C21").Formula = "=SUM(C18:C20)"

$ws.Range("C22").Value = -33376095

$ws.Range("C26").Value = 1003368421

$wb.Save()
